# Updated symbol list (coinranking snapshot refresh).
# For each changed cell we re-assign the new scraped value. Numeric-looking
# text (prices / percentages) is written with NumberFormat "@" set first so
# Excel keeps storing it as text (e.g. "41.50", "1,897.31%") instead of
# silently coercing it to a Double and losing the original formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "319.63"
$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = "3.31%"
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "41.50"
$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = "1.42%"
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "0.07745"
$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = "1.49%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "4.358"
$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = "1.51%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.742"
$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = "8.20%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.9453"
$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = "4.04%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "2.425"
$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = "-1.89%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.1242"
$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = "-2.31%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "0.1865"
$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = "3.46%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.09245"
$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = "2.40%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.04107"
$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = "-5.49%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.1053"
$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = "0.60%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.001282"
$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = "2.51%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "0.005751"
$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = "1.71%"
$ws.Range("B17").Value = "UpBots"
$ws.Range("C17").Value = "https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt"
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "0.007491"
$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = "1,897.31%"
$ws.Range("B18").Value = "LEO"
$ws.Range("C18").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.352"
$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = "-0.04%"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "0.3358"
$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = "1.33%"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "8.395"
$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = "21.49%"
$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = "-2.31%"
$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = "3.03%"
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.04024"
$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = "-0.43%"
$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = "-0.16%"
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "0.004132"
$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = "2.01%"
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "0.0001272"
$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = "0.13%"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.02562"
$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = "5.94%"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.05350"
$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = "2.44%"
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.007790"
$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = "-0.69%"
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.1318"
$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = "1.29%"
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.007044"
$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = "3.62%"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.001991"
$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = "8.24%"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "0.008252"
$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = "11.14%"
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.3169"
$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = "-5.44%"
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "0.00006699"
$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = "-2.58%"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.00000000751"
$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = "-0.03%"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.1996"
$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = "25.28%"
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.004206"
$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = "40.07%"
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "0.00002103"
$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = "-0.03%"
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.0002003"
$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = "-0.03%"
